$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column to remain text so numeric-looking values
# like "604.03" are not auto-converted to floating point numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "63.950.37"
$ws.Range("E2").Value = "  -1.28%  "
$ws.Range("D3").Value = "3.149.65"
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "604.03"
$ws.Range("E5").Value = "  -2.00%  "
$ws.Range("D6").Value = "143.52"
$ws.Range("E6").Value = "  -2.48%  "
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("D8").Value = "3.144.63"
$ws.Range("E8").Value = "  -0.55%  "
$ws.Range("D9").Value = "0.527"
$ws.Range("E9").Value = "  -0.47%  "
$ws.Range("E10").Value = "  -1.83%  "
$ws.Range("D11").Value = "5.39"
$ws.Range("E11").Value = "  -1.98%  "
$ws.Range("D12").Value = "0.467"
$ws.Range("E12").Value = "  -1.61%  "
$ws.Range("E13").Value = "  -2.27%  "
$ws.Range("E14").Value = "  -2.37%  "
$ws.Range("D15").Value = "3.664.92"
$ws.Range("E15").Value = "  -0.57%  "
$ws.Range("E16").Value = "  +2.62%  "
$ws.Range("D17").Value = "63.995.03"
$ws.Range("E17").Value = "  -1.18%  "
$ws.Range("D18").Value = "3.158.98"
$ws.Range("E18").Value = "  -0.07%  "
$ws.Range("D19").Value = "6.87"
$ws.Range("E19").Value = "  -0.95%  "
$ws.Range("D20").Value = "489.67"
$ws.Range("E21").Value = "  -0.10%  "
$ws.Range("E22").Value = "  -1.91%  "
$ws.Range("E23").Value = "  -3.27%  "
$ws.Range("D24").Value = "88.11"
$ws.Range("E24").Value = "  +4.31%  "
$ws.Range("E25").Value = "  -3.61%  "
$ws.Range("E26").Value = "  -0.01%  "
$ws.Range("E27").Value = "  -2.43%  "
$ws.Range("E28").Value = "  -4.22%  "
$ws.Range("D29").Value = "6.97"
$ws.Range("E29").Value = "  +0.71%  "
$ws.Range("E30").Value = "  -1.24%  "
$ws.Range("D31").Value = "27.78"
$ws.Range("E31").Value = "  +4.58%  "
$ws.Range("E32").Value = "  -5.65%  "
$ws.Range("E33").Value = "  -0.01%  "
$ws.Range("E34").Value = "  -1.97%  "
$ws.Range("E35").Value = "  -2.85%  "
$ws.Range("E36").Value = "  +0.60%  "
$ws.Range("D37").Value = "52.67"
$ws.Range("E37").Value = "  -0.73%  "
$ws.Range("D38").Value = "0.0₃0739"
$ws.Range("E38").Value = "  -5.07%  "
$ws.Range("E39").Value = "  -6.70%  "
$ws.Range("D40").Value = "434.41"
$ws.Range("E40").Value = "  -5.62%  "
$ws.Range("D41").Value = "0.0397"
$ws.Range("E41").Value = "  -0.61%  "
$ws.Range("E42").Value = "  -0.44%  "
$ws.Range("D43").Value = "8.32"
$ws.Range("E43").Value = "  -1.08%  "
$ws.Range("D44").Value = "2.935.76"
$ws.Range("E44").Value = "  +3.04%  "
$ws.Range("E45").Value = "  -3.33%  "
$ws.Range("E46").Value = "  -5.50%  "
$ws.Range("D47").Value = "2.40"
$ws.Range("E47").Value = "  -1.39%  "
$ws.Range("D49").Value = "25.85"
$ws.Range("E49").Value = "  -2.81%  "
$ws.Range("E50").Value = "  -0.12%  "
$ws.Range("D51").Value = "120.30"
$ws.Range("E51").Value = "  -0.01%  "
